$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 49
$wsExhibit.Range("F4").Value = 3542
$wsExhibit.Range("F5").Value = 2202
$wsExhibit.Range("F9").Value = 60
$wsExhibit.Range("F12").Value = 1805

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 49
$wsAll.Range("F4").Value = 3542
$wsAll.Range("F5").Value = 2202
$wsAll.Range("F10").Value = 60
$wsAll.Range("F15").Value = 1805
